# Fix typo "wondows" -> "windows" in the "Added this line in wondows" paragraph.
# The target edit (per the diff) keeps the text as separate runs split right
# where the correction happened, with the "_GoBack" bookmark repositioned to
# mark that edit point (mirroring Word's own "last edit" bookmark behavior).

$d = $word.ActiveDocument

# Locate the paragraph that needs the fix.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*wondows*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find paragraph containing 'wondows'"
}

$pStart = $target.Range.Start

# Offset of the "o" that must become "i" ("Added this line in w" is 20 chars).
$oOffset = $pStart + 20

# Replace just that one character, with revision tracking on so the edit
# naturally splits the surrounding text into separate runs instead of
# collapsing everything back into one merged run.
$d.TrackRevisions = $true
$charRng = $d.Range($oOffset, $oOffset + 1)
$charRng.Text = "i"
$d.TrackRevisions = $false
$d.Revisions.AcceptAll()

# Move the "_GoBack" bookmark to sit right after the inserted "i" (matching
# Word's behavior of leaving _GoBack at the most recent edit location).
$editPoint = $d.Range($oOffset + 1, $oOffset + 1)
$d.Bookmarks.Add("_GoBack", $editPoint)
